# Auto-generated edit script: update cached market-price values on the
# Leviathan_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 5068
$ws.Range("I98").Value = 3932.8125
$ws.Range("K98").Value = 3932.8125
$ws.Range("M98").Value = -2434.8125
$ws.Range("H122").Value = 5068
$ws.Range("I122").Value = 3932.8125
$ws.Range("K122").Value = 11798.4375
$ws.Range("M122").Value = -9348.4375
$ws.Range("H132").Value = 4270.484
$ws.Range("I132").Value = 1533
$ws.Range("K132").Value = 4599
$ws.Range("M132").Value = -2069
$ws.Range("H138").Value = 2222.5
$ws.Range("I138").Value = 1444.3846
$ws.Range("J138").Value = 3346.4443
$ws.Range("K138").Value = 4333.1538
$ws.Range("L138").Value = 10039.3329
$ws.Range("M138").Value = 806.8462
$ws.Range("N138").Value = -20319.3329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H61").Value = 3635.1
$ws.Range("I61").Value = 3682.111
$ws.Range("J61").Value = 3212
$ws.Range("K61").Value = 3682.111
$ws.Range("L61").Value = 3212
$ws.Range("M61").Value = -3470.111
$ws.Range("N61").Value = -3636
$ws.Range("H74").Value = 1395.1471
$ws.Range("I74").Value = 1114.5
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 1114.5
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -240.5
$ws.Range("N74").Value = -5248
$ws.Range("H77").Value = 1395.1471
$ws.Range("I77").Value = 1114.5
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 5572.5
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -1204.5
$ws.Range("N77").Value = -26236
$ws.Range("H122").Value = 11321.19
$ws.Range("I122").Value = 12921.056
$ws.Range("K122").Value = 38763.16800000001
$ws.Range("M122").Value = -36313.16800000001
$ws.Range("H132").Value = 1473.3334
$ws.Range("I132").Value = 964.65216
$ws.Range("K132").Value = 2893.95648
$ws.Range("M132").Value = -363.9564799999998
$ws.Range("H136").Value = 3635.1
$ws.Range("I136").Value = 3682.111
$ws.Range("J136").Value = 3212
$ws.Range("K136").Value = 11046.333
$ws.Range("L136").Value = 9636
$ws.Range("M136").Value = -8496.332999999999
$ws.Range("N136").Value = -14736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8511.741
$ws.Range("I20").Value = 6350.2354
$ws.Range("J20").Value = 12186.3
$ws.Range("K20").Value = 6350.2354
$ws.Range("L20").Value = 12186.3
$ws.Range("M20").Value = -6103.2354
$ws.Range("N20").Value = -12680.3
$ws.Range("H86").Value = 1930.6
$ws.Range("I86").Value = 1687.375
$ws.Range("K86").Value = 1687.375
$ws.Range("M86").Value = -564.375
$ws.Range("H89").Value = 1930.6
$ws.Range("I89").Value = 1687.375
$ws.Range("K89").Value = 8436.875
$ws.Range("M89").Value = -2820.875
$ws.Range("H105").Value = 5004327.5
$ws.Range("J105").Value = 2580.25
$ws.Range("L105").Value = 2580.25
$ws.Range("N105").Value = -6074.25
$ws.Range("H134").Value = 2521
$ws.Range("I134").Value = 2521
$ws.Range("K134").Value = 7563
$ws.Range("M134").Value = -5028

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 203.47826
$ws.Range("I22").Value = 198.94737
$ws.Range("J22").Value = 225
$ws.Range("K22").Value = 198.94737
$ws.Range("L22").Value = 225
$ws.Range("M22").Value = 151.05263
$ws.Range("N22").Value = -925

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 13996
$ws.Range("J62").Value = 14994.5
$ws.Range("L62").Value = 44983.5
$ws.Range("N62").Value = -46355.5
$ws.Range("H65").Value = 13996
$ws.Range("J65").Value = 14994.5
$ws.Range("L65").Value = 134950.5
$ws.Range("N65").Value = -141814.5
$ws.Range("H121").Value = 20835438
$ws.Range("J121").Value = 3204
$ws.Range("L121").Value = 9612
$ws.Range("N121").Value = -12232

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 12090.444
$ws.Range("I99").Value = 9786.875
$ws.Range("K99").Value = 9786.875
$ws.Range("M99").Value = -7540.875
$ws.Range("H102").Value = 2801.1333
$ws.Range("I102").Value = 2956.0908
$ws.Range("K102").Value = 2956.0908
$ws.Range("M102").Value = -1334.0908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 40699.715
$ws.Range("I7").Value = 54780
$ws.Range("K7").Value = 54780
$ws.Range("M7").Value = -54668
$ws.Range("H22").Value = 1372.5294
$ws.Range("J22").Value = 1159.1111
$ws.Range("L22").Value = 1159.1111
$ws.Range("N22").Value = -1749.1111
$ws.Range("H23").Value = 20766.5
$ws.Range("I23").Value = 17919.8
$ws.Range("K23").Value = 17919.8
$ws.Range("M23").Value = -17689.8
$ws.Range("H27").Value = 1372.5294
$ws.Range("J27").Value = 1159.1111
$ws.Range("L27").Value = 1159.1111
$ws.Range("N27").Value = -1373.1111
$ws.Range("H40").Value = 4861.5454
$ws.Range("I40").Value = 4165.4443
$ws.Range("J40").Value = 7994
$ws.Range("K40").Value = 4165.4443
$ws.Range("L40").Value = 7994
$ws.Range("M40").Value = -4029.4443
$ws.Range("N40").Value = -8266
$ws.Range("H46").Value = 36696.668
$ws.Range("I46").Value = 53989
$ws.Range("K46").Value = 53989
$ws.Range("M46").Value = -53801
$ws.Range("H55").Value = 301.2
$ws.Range("J55").Value = 403.1875
$ws.Range("L55").Value = 403.1875
$ws.Range("N55").Value = -749.1875
$ws.Range("H61").Value = 30337542
$ws.Range("I61").Value = 37041516
$ws.Range("K61").Value = 37041516
$ws.Range("M61").Value = -37041314
$ws.Range("H68").Value = 2383.9473
$ws.Range("I68").Value = 2098.8667
$ws.Range("K68").Value = 2098.8667
$ws.Range("M68").Value = -1349.8667
$ws.Range("H71").Value = 2383.9473
$ws.Range("I71").Value = 2098.8667
$ws.Range("K71").Value = 10494.3335
$ws.Range("M71").Value = -6750.333500000001
$ws.Range("H113").Value = 30337542
$ws.Range("I113").Value = 37041516
$ws.Range("K113").Value = 37041516
$ws.Range("M113").Value = -37039346
$ws.Range("H122").Value = 12540.3
$ws.Range("I122").Value = 17317.166
$ws.Range("J122").Value = 5375
$ws.Range("K122").Value = 51951.49800000001
$ws.Range("L122").Value = 16125
$ws.Range("M122").Value = -49501.49800000001
$ws.Range("N122").Value = -21025
$ws.Range("H126").Value = 40699.715
$ws.Range("I126").Value = 54780
$ws.Range("K126").Value = 164340
$ws.Range("M126").Value = -161870
$ws.Range("H132").Value = 6527.0386
$ws.Range("I132").Value = 6661.2173
$ws.Range("J132").Value = 5498.3335
$ws.Range("K132").Value = 19983.6519
$ws.Range("L132").Value = 16495.0005
$ws.Range("M132").Value = -17453.6519
$ws.Range("N132").Value = -21555.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H81").Value = 2249.3333
$ws.Range("J81").Value = 3065.3333
$ws.Range("L81").Value = 6130.6666
$ws.Range("N81").Value = -8252.6666
$ws.Range("H84").Value = 2249.3333
$ws.Range("J84").Value = 3065.3333
$ws.Range("L84").Value = 30653.333
$ws.Range("N84").Value = -41261.333
$ws.Range("H113").Value = 1373.3636
$ws.Range("I113").Value = 1230
$ws.Range("K113").Value = 3690
$ws.Range("M113").Value = -1520
$ws.Range("H122").Value = 6219.1113
$ws.Range("I122").Value = 6371.5
$ws.Range("K122").Value = 19114.5
$ws.Range("M122").Value = -16664.5
$ws.Range("H132").Value = 7893.8335
$ws.Range("I132").Value = 8011.1177
$ws.Range("K132").Value = 24033.3531
$ws.Range("M132").Value = -21503.3531
$ws.Range("H136").Value = 2102.1614
$ws.Range("I136").Value = 1770.72
$ws.Range("K136").Value = 5312.16
$ws.Range("M136").Value = -2762.16
